$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Mapping of row -> column -> new value, derived from the diff
$data = @{
    2 = @{ "B" = 19.95221006506002; "C" = 9.523817997876336; "D" = 7.263747435841142; "E" = 9.38983413409801; "F" = 39.11008160143938; "I" = 30.67827950301434; "L" = 10.49374847646193 }
    3 = @{ "B" = 19.51527105547694; "C" = 8.959559772404363; "D" = 7.282965234280929; "E" = 9.405329733874636; "F" = 38.73492891471737; "I" = 30.58314328341859; "L" = 10.47525806446442 }
    4 = @{ "B" = 19.2485225153784; "C" = 8.596659747759153; "D" = 7.295720102052071; "E" = 9.415484258095381; "F" = 38.51459485985152; "I" = 30.53209033613397; "L" = 10.46620087574238 }
    5 = @{ "B" = 19.14038131367918; "C" = 8.444775081115507; "D" = 7.301156872693726; "E" = 9.419783543343966; "F" = 38.42740137791335; "I" = 30.51313955152014; "L" = 10.46308911732526 }
    6 = @{ "B" = 19.1224636123794; "C" = 8.419317679726573; "D" = 7.302074050889313; "E" = 9.420507182230903; "F" = 38.41308179345099; "I" = 30.51010478251471; "L" = 10.46260743214769 }
    7 = @{ "B" = 19.24706158014774; "C" = 8.594627374023107; "D" = 7.295792457780713; "E" = 9.415541586568741; "F" = 38.51340833639536; "I" = 30.53182725196275; "L" = 10.46615656266563 }
    8 = @{ "B" = 19.80133929102289; "C" = 9.332746818840858; "D" = 7.270174846210807; "E" = 9.39504430710369; "F" = 38.9786991849553; "I" = 30.64395101892245; "L" = 10.48689764756647 }
    9 = @{ "B" = 20.89294613006399; "C" = 10.7030257932566; "D" = 7.227566252419821; "E" = 9.359916952761653; "F" = 39.96690302714735; "I" = 30.92201830630472; "L" = 10.54569146654218 }
    10 = @{ "B" = 21.68815800892214; "C" = 11.65592536533999; "D" = 7.20097994422631; "E" = 9.337182297514634; "F" = 40.73387671377539; "I" = 31.16136813160038; "L" = 10.59977624191555 }
    11 = @{ "B" = 22.04655271727577; "C" = 12.06469726788548; "D" = 7.189924629262086; "E" = 9.327503741308254; "F" = 41.09048443099115; "I" = 31.27773105724925; "L" = 10.62670200037054 }
    12 = @{ "B" = 22.18163975910063; "C" = 12.21595684453951; "D" = 7.185888941332712; "E" = 9.323933896757771; "F" = 41.2265262467049; "I" = 31.3228555317661; "L" = 10.63722759764451 }
    13 = @{ "B" = 22.15257657828244; "C" = 12.18353723291766; "D" = 7.186751374371801; "E" = 9.324698495151502; "F" = 41.19718433960495; "I" = 31.31309027665803; "L" = 10.63494614906975 }
    14 = @{ "B" = 22.05767983074782; "C" = 12.07721218969031; "D" = 7.189589582499149; "E" = 9.327208140953418; "F" = 41.10165720053386; "I" = 31.28142231951274; "L" = 10.62756137361983 }
    15 = @{ "B" = 21.99946667488283; "C" = 12.01162544166569; "D" = 7.191347734661573; "E" = 9.328757765272181; "F" = 41.04327141083554; "I" = 31.26216237793215; "L" = 10.623080740036 }
    16 = @{ "B" = 21.66465582365627; "C" = 11.62871564968213; "D" = 7.201723453697467; "E" = 9.337828148803625; "F" = 40.710717623063; "I" = 31.15391281154296; "L" = 10.59806293253555 }
    17 = @{ "B" = 21.45829404364759; "C" = 11.38749865246189; "D" = 7.208355652137945; "E" = 9.343562340819791; "F" = 40.50860579131517; "I" = 31.0894110407357; "L" = 10.58330702486917 }
    18 = @{ "B" = 21.33929300318195; "C" = 11.24643033151869; "D" = 7.212268026150713; "E" = 9.346922964989735; "F" = 40.39308836884224; "I" = 31.05301658005054; "L" = 10.57503870625499 }
    19 = @{ "B" = 21.29895317885869; "C" = 11.19826682144453; "D" = 7.213609430306847; "E" = 9.348071549903056; "F" = 40.35410507375844; "I" = 31.04081560647256; "L" = 10.57227692139444 }
    20 = @{ "B" = 21.48029449949723; "C" = 11.41341727615806; "D" = 7.207639521266566; "E" = 9.342945462518736; "F" = 40.53004594952814; "I" = 31.09620447908702; "L" = 10.58485519378523 }
    21 = @{ "B" = 22.085571510544; "C" = 12.10853817571903; "D" = 7.18875183084801; "E" = 9.326468414647401; "F" = 41.12968949046338; "I" = 31.29069532156327; "L" = 10.6297215598927 }
    22 = @{ "B" = 22.47742571362603; "C" = 12.54225302758599; "D" = 7.177286770266772; "E" = 9.316254604060301; "F" = 41.52738273510404; "I" = 31.42397834378262; "L" = 10.66096188282128 }
    23 = @{ "B" = 22.26867317683481; "C" = 12.31264856866929; "D" = 7.183325006219405; "E" = 9.321655200259265; "F" = 41.31463185786316; "I" = 31.352283572781; "L" = 10.64411449410316 }
    24 = @{ "B" = 21.4703492062582; "C" = 11.401706910473; "D" = 7.207962974416557; "E" = 9.343224153959442; "F" = 40.52035072990979; "I" = 31.0931310200398; "L" = 10.584154596448 }
    25 = @{ "B" = 20.59815977354552; "C" = 10.33007424475665; "D" = 7.238269434440639; "E" = 9.368878968164234; "F" = 39.69198121533708; "I" = 30.84059858920527; "L" = 10.52786066065584 }
}

foreach ($row in $data.Keys) {
    $rowData = $data[$row]
    foreach ($col in $rowData.Keys) {
        $ws.Range("$col$row").Value = $rowData[$col]
    }
}
